# Documento Release – Proyecto APT
# Actualización de Sprints Semanales:
#  - Release 1.2: la fecha estimada se extiende del 15 al 24 de octubre,
#    y el estado pasa de "Completado" a "En curso".
#  - Release 1.5: ya no se usa Trello para el registro, solo GitHub.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Release 1.2 – Diseño y Arquitectura de Software (fila 9)
$ws.Range("C9").Value = "22 de septiembre – 24 de octubre 2025"
$ws.Range("F9").Value = "En curso"

# Release 1.5 – Gestión de Proyecto y Documentación Continua (fila 12)
$ws.Range("D12").Value = "- Registro constante en GitHub.`n- Seguimiento de Sprints, reuniones y control de versiones.`n- Documentación de cambios y evidencias."
